# edit.ps1
#
# Recreates the template-authoring edit that added six new demo sheets
# ("No PAE", "PAE Clear", "PAE Remove", "groupDirNone", "groupDirRows",
# "groupDirCols") to ImplCollProcessingTemplate.xlsx, each showing a
# variant of JETT's implicit-collection-processing / grouping syntax.

$wb = $excel.ActiveWorkbook

# Colors (Excel COM uses 0xBBGGRR ordering for RGB(r,g,b))
$BLUE = 153 + 204*256 + 255*65536      # #99CCFF
$GRAY = 192 + 192*256 + 192*65536      # #C0C0C0

function Set-TitleCell($cell) {
    $cell.Borders.Color = 0
    $cell.Borders.Weight = 2
    $cell.Borders.LineStyle = 1
    $cell.Interior.Color = $BLUE
    $cell.Font.Bold = $true
}

function Set-HeaderCell($cell) {
    $cell.Borders.Color = 0
    $cell.Borders.Weight = 2
    $cell.Borders.LineStyle = 1
    $cell.Interior.Color = $BLUE
    $cell.Font.Bold = $true
}

function Set-DataCell($cell) {
    $cell.Borders.Color = 0
    $cell.Borders.Weight = 2
    $cell.Borders.LineStyle = 1
    $cell.Interior.Color = $GRAY
}

function Set-DataNumericCell($cell) {
    $cell.Borders.Color = 0
    $cell.Borders.Weight = 2
    $cell.Borders.LineStyle = 1
    $cell.Interior.Color = $GRAY
    $cell.NumberFormat = "0.000"
}

# Style + fill in a standard 5-column "City/Name/Wins/Losses/Pct." block
# whose top-left corner is ($row, $col) on worksheet $ws :
#   row      : merged title cell   (s=4-equivalent)
#   row + 1  : header cells        (s=1-equivalent)
#   row + 2  : data cells          (s=2 / s=3-equivalent), values from $dataVals
function Fill-Block($ws, $row, $col, $titleText, $dataVals) {
    # Title row (merged across 5 columns), styled cell-by-cell before merging
    for ($c = $col; $c -le ($col + 4); $c++) {
        Set-TitleCell($ws.Cells.Item($row, $c))
    }
    $ws.Cells.Item($row, $col).Value = $titleText

    # Header row
    $headers = @("City", "Name", "Wins", "Losses", "Pct.")
    for ($i = 0; $i -lt 5; $i++) {
        $cell = $ws.Cells.Item($row + 1, $col + $i)
        $cell.Value = $headers[$i]
        Set-HeaderCell($cell)
    }

    # Data row
    for ($i = 0; $i -lt 5; $i++) {
        $cell = $ws.Cells.Item($row + 2, $col + $i)
        $cell.Value = $dataVals[$i]
        if ($i -eq 4) {
            Set-DataNumericCell($cell)
        } else {
            Set-DataCell($cell)
        }
    }

    # Merge title row last, after styling is applied to each cell
    $ws.Range($ws.Cells.Item($row, $col), $ws.Cells.Item($row, $col + 4)).Merge()
}

function Set-ColWidth($ws, $colIndex, $targetRawWidth) {
    # This runtime stores column width (raw OOXML "width" attribute) as
    # round(ColumnWidth * 6 + 5) / 6 ; invert that as closely as possible.
    $chars = ($targetRawWidth * 6 - 5) / 6
    $ws.Columns.Item($colIndex).ColumnWidth = $chars
}

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------------------
# Sheet: "No PAE" -- two independent implicit-collection loops side by side,
# neither one using @pastEndAction.
# ---------------------------------------------------------------------------
$wsNoPae = $wb.Worksheets.Add($null, $lastSheet)
$wsNoPae.Name = "No PAE"
Fill-Block $wsNoPae 1 1 'Division: ${ofTheirOwn.name}' @('${ofTheirOwn.teams.city}', '${ofTheirOwn.teams.name}', '${ofTheirOwn.teams.wins}', '${ofTheirOwn.teams.losses}', '${ofTheirOwn.teams.pct}')
Fill-Block $wsNoPae 1 6 'Division: ${pacific.name}' @('${pacific.teams.city}', '${pacific.teams.name}', '${pacific.teams.wins}', '${pacific.teams.losses}', '${pacific.teams.pct}')
Set-ColWidth $wsNoPae 1 15
Set-ColWidth $wsNoPae 2 15.5703125
Set-ColWidth $wsNoPae 6 15
Set-ColWidth $wsNoPae 7 15.5703125
[void]$wsNoPae.Range("A1:E1").Select()
$lastSheet = $wsNoPae

# ---------------------------------------------------------------------------
# Sheet: "PAE Clear" -- same as above, but left loop uses
# @pastEndAction=clear
# ---------------------------------------------------------------------------
$wsPaeClear = $wb.Worksheets.Add($null, $lastSheet)
$wsPaeClear.Name = "PAE Clear"
Fill-Block $wsPaeClear 1 1 'Division: ${ofTheirOwn.name}' @('${ofTheirOwn.teams.city}?@pastEndAction=clear', '${ofTheirOwn.teams.name}', '${ofTheirOwn.teams.wins}', '${ofTheirOwn.teams.losses}', '${ofTheirOwn.teams.pct}')
Fill-Block $wsPaeClear 1 6 'Division: ${pacific.name}' @('${pacific.teams.city}', '${pacific.teams.name}', '${pacific.teams.wins}', '${pacific.teams.losses}', '${pacific.teams.pct}')
Set-ColWidth $wsPaeClear 1 15
Set-ColWidth $wsPaeClear 2 15.5703125
Set-ColWidth $wsPaeClear 6 15
Set-ColWidth $wsPaeClear 7 15.5703125
[void]$wsPaeClear.Range("A1:E1").Select()
$lastSheet = $wsPaeClear

# ---------------------------------------------------------------------------
# Sheet: "PAE Remove" -- same as above, but left loop uses
# @pastEndAction=remove
# ---------------------------------------------------------------------------
$wsPaeRemove = $wb.Worksheets.Add($null, $lastSheet)
$wsPaeRemove.Name = "PAE Remove"
Fill-Block $wsPaeRemove 1 1 'Division: ${ofTheirOwn.name}' @('${ofTheirOwn.teams.city}?@pastEndAction=remove', '${ofTheirOwn.teams.name}', '${ofTheirOwn.teams.wins}', '${ofTheirOwn.teams.losses}', '${ofTheirOwn.teams.pct}')
Fill-Block $wsPaeRemove 1 6 'Division: ${pacific.name}' @('${pacific.teams.city}', '${pacific.teams.name}', '${pacific.teams.wins}', '${pacific.teams.losses}', '${pacific.teams.pct}')
Set-ColWidth $wsPaeRemove 1 15
Set-ColWidth $wsPaeRemove 2 15.5703125
Set-ColWidth $wsPaeRemove 6 15
Set-ColWidth $wsPaeRemove 7 15.5703125
[void]$wsPaeRemove.Range("A1:E1").Select()
$lastSheet = $wsPaeRemove

# ---------------------------------------------------------------------------
# Sheet: "groupDirNone" -- single implicit-collection loop using
# @groupDir=none;collapse=...
# ---------------------------------------------------------------------------
$wsGroupNone = $wb.Worksheets.Add($null, $lastSheet)
$wsGroupNone.Name = "groupDirNone"
Fill-Block $wsGroupNone 1 1 'Division: ${divisionsList.name}?@extraRows=2' @("`${divisionsList.teams.city}?@groupDir=none;collapse=`${divisionsList.name.equals('Southeast')}", '${divisionsList.teams.name}', '${divisionsList.teams.wins}', '${divisionsList.teams.losses}', '${divisionsList.teams.pct}')
Set-ColWidth $wsGroupNone 1 15
Set-ColWidth $wsGroupNone 2 15.5703125
[void]$wsGroupNone.Range("A1:E1").Select()
$lastSheet = $wsGroupNone

# ---------------------------------------------------------------------------
# Sheet: "groupDirRows" -- single implicit-collection loop using
# @groupDir=rows;collapse=...
# ---------------------------------------------------------------------------
$wsGroupRows = $wb.Worksheets.Add($null, $lastSheet)
$wsGroupRows.Name = "groupDirRows"
Fill-Block $wsGroupRows 1 1 'Division: ${divisionsList.name}?@extraRows=2' @("`${divisionsList.teams.city}?@groupDir=rows;collapse=`${divisionsList.name.equals('Southeast')}", '${divisionsList.teams.name}', '${divisionsList.teams.wins}', '${divisionsList.teams.losses}', '${divisionsList.teams.pct}')
Set-ColWidth $wsGroupRows 1 15
Set-ColWidth $wsGroupRows 2 15.5703125
[void]$wsGroupRows.Range("A1:E1").Select()
$lastSheet = $wsGroupRows

# ---------------------------------------------------------------------------
# Sheet: "groupDirCols" -- single implicit-collection loop using
# @groupDir=cols;collapse=..., with the block offset one column to the
# right (starts at column B) and an extra narrow column A.
# ---------------------------------------------------------------------------
$wsGroupCols = $wb.Worksheets.Add($null, $lastSheet)
$wsGroupCols.Name = "groupDirCols"
Fill-Block $wsGroupCols 1 2 'Division: ${divisionsList.name}?@extraRows=2;left=1;right=4;copyRight=true' @("`${divisionsList.teams.city}?@left=0;right=4;groupDir=cols;collapse=`${divisionsList.name.equals('Southeast')}", '${divisionsList.teams.name}', '${divisionsList.teams.wins}', '${divisionsList.teams.losses}', '${divisionsList.teams.pct}')
Set-ColWidth $wsGroupCols 1 4.140625
Set-ColWidth $wsGroupCols 2 15
Set-ColWidth $wsGroupCols 3 15.5703125
$lastSheet = $wsGroupCols

# Restore the original active sheet/tab so the workbook opens the same way
# it did before this edit.
[void]$wb.Worksheets.Item("Implicit").Select()
